# Fill in the new "v1282" performance column (F) on the Sponza and
# ComplexMesh sheets, which previously existed only as an empty
# placeholder column (header + formulas already present, but blank
# data caused #DIV/0! errors).

$wb = $excel.ActiveWorkbook

$sponza = $wb.Worksheets.Item("Sponza")
$complexMesh = $wb.Worksheets.Item("ComplexMesh")

# --- Sponza ("PipelinePerformance" sheet 2) ---------------------------
$sponza.Range("F1").Value = "v1282"

$sponzaValues = 10031, 10058, 10015, 10048, 10141, 10045, 10021, 10039, 10043, 10041
for ($i = 0; $i -lt $sponzaValues.Length; $i++) {
    $sponza.Cells.Item(2 + $i, 6).Value = $sponzaValues[$i]
}

# Step-over-step comparison now walks E->F instead of D->F now that F
# is populated.
$sponza.Range("F15").Formula = "=E12/F12"

# --- ComplexMesh (sheet 3) ---------------------------------------------
$complexMesh.Range("F1").Value = "v1282"

$complexMeshValues = 7530, 7551, 7535, 7513, 7517, 7538, 7553, 7544, 7586, 7573
for ($i = 0; $i -lt $complexMeshValues.Length; $i++) {
    $complexMesh.Cells.Item(2 + $i, 6).Value = $complexMeshValues[$i]
}

$complexMesh.Range("F15").Formula = "=E12/F12"

# --- Navigation: leave a selection mark on Sponza, then switch focus
#     over to ComplexMesh (now the active tab) with its own selection.
$sponza.Range("D3").Select() | Out-Null

$complexMesh.Activate() | Out-Null
$complexMesh.Range("F15").Select() | Out-Null
